# Time Sheet.xlsx update:
#  - Split the single "Evan" person into "Evan F" and "Evan S"
#  - Re-label the existing first time-sheet entry (row 3) as "Evan S"
#  - Add a new time-sheet entry (row 4) for "Evan F" with a new task
#  - Move the active cell selection to E5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Create the "Evan F" entry first so the new shared string is registered
# before "Evan S" (keeps shared-string ordering in line with the source).
$ws.Range("A4").Value = "Evan F"

# The person on the original first entry is actually "Evan S".
$ws.Range("A3").Value = "Evan S"

# New row 4: date worked, start/end times and the task description.
$ws.Range("B4").NumberFormat = "mm-dd-yy"
$ws.Range("B4").Value = 42771

$ws.Range("C4").NumberFormat = "h:mm"
$ws.Range("C4").Value = 0.4375

$ws.Range("D4").NumberFormat = "h:mm"
$ws.Range("D4").Value = 0.45833333333333331

$ws.Range("E4").Value = "Made smoother green grass and made chest sprites"

# Match the saved cursor position from the source file.
$ws.Range("E5").Select()
